# KHL stats runtime update (2025-10-28 11:08:25)
# - Append 3 new matches (rows 398-400) to Matches_SOG with their shots-on-goal.
# - Refresh as_of_utc timestamps on Shots_HA / Shots_Summary / Meta_ext, and
#   recompute the aggregate stats for the teams involved in the new matches.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Matches_SOG: append the three newly played games.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Matches_SOG")

$newMatches = @(
    @{Row=398; Uid="897696"; Date="2025-10-27T16:30:00"; Home="Авангард";   Away="Металлург Мг"; SogHome=32; SogAway=28},
    @{Row=399; Uid="897695"; Date="2025-10-27T19:00:00"; Home="Нефтехимик"; Away="Барыс";         SogHome=34; SogAway=32},
    @{Row=400; Uid="897697"; Date="2025-10-27T19:00:00"; Home="Ак Барс";    Away="Адмирал";       SogHome=25; SogAway=37}
)

foreach ($m in $newMatches) {
    $ws.Cells.Item($m.Row, 1).Value = "'" + $m.Uid
    $ws.Cells.Item($m.Row, 2).Value = $m.Date
    $ws.Cells.Item($m.Row, 3).Value = $m.Home
    $ws.Cells.Item($m.Row, 4).Value = $m.Away
    $ws.Cells.Item($m.Row, 5).Value = $m.SogHome
    $ws.Cells.Item($m.Row, 6).Value = $m.SogAway
    $ws.Cells.Item($m.Row, 7).Value = "khl_text"
}

# ---------------------------------------------------------------------------
# 2) Shots_HA: bump as_of_utc for every team, and update the home/away shot
#    totals for the six teams that played the new games.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Shots_HA")

for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 4).Value = "2025-10-27T19:00:00Z"
}

$ws.Cells.Item(2, 5).Value = 19
$ws.Cells.Item(2, 7).Value = 636
$ws.Cells.Item(2, 8).Value = 522
$ws.Cells.Item(2, 9).Value = 33.5
$ws.Cells.Item(2, 10).Value = 27.5

$ws.Cells.Item(4, 6).Value = 17
$ws.Cells.Item(4, 11).Value = 537
$ws.Cells.Item(4, 12).Value = 499
$ws.Cells.Item(4, 13).Value = 31.6
$ws.Cells.Item(4, 14).Value = 29.4

$ws.Cells.Item(5, 5).Value = 21
$ws.Cells.Item(5, 7).Value = 711
$ws.Cells.Item(5, 8).Value = 533
$ws.Cells.Item(5, 9).Value = 33.9
$ws.Cells.Item(5, 10).Value = 25.4

$ws.Cells.Item(7, 6).Value = 13
$ws.Cells.Item(7, 11).Value = 350
$ws.Cells.Item(7, 12).Value = 452
$ws.Cells.Item(7, 13).Value = 26.9

$ws.Cells.Item(13, 6).Value = 13
$ws.Cells.Item(13, 11).Value = 372
$ws.Cells.Item(13, 12).Value = 336
$ws.Cells.Item(13, 13).Value = 28.6
$ws.Cells.Item(13, 14).Value = 25.8

$ws.Cells.Item(14, 5).Value = 23
$ws.Cells.Item(14, 7).Value = 724
$ws.Cells.Item(14, 8).Value = 770
$ws.Cells.Item(14, 9).Value = 31.5

# ---------------------------------------------------------------------------
# 3) Shots_Summary: bump as_of_utc for every team, and update the combined
#    shot totals for the six teams that played the new games.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Shots_Summary")

for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 4).Value = "2025-10-27T19:00:00Z"
}

$ws.Cells.Item(2, 5).Value = 35
$ws.Cells.Item(2, 6).Value = 1208
$ws.Cells.Item(2, 7).Value = 994
$ws.Cells.Item(2, 8).Value = 34.5

$ws.Cells.Item(4, 5).Value = 31
$ws.Cells.Item(4, 6).Value = 1081
$ws.Cells.Item(4, 7).Value = 891
$ws.Cells.Item(4, 8).Value = 34.9
$ws.Cells.Item(4, 9).Value = 28.7

$ws.Cells.Item(5, 5).Value = 39
$ws.Cells.Item(5, 6).Value = 1315
$ws.Cells.Item(5, 7).Value = 1067
$ws.Cells.Item(5, 8).Value = 33.7
$ws.Cells.Item(5, 9).Value = 27.4

$ws.Cells.Item(7, 5).Value = 39
$ws.Cells.Item(7, 6).Value = 1162
$ws.Cells.Item(7, 7).Value = 1278
$ws.Cells.Item(7, 8).Value = 29.8
$ws.Cells.Item(7, 9).Value = 32.8

$ws.Cells.Item(13, 5).Value = 37
$ws.Cells.Item(13, 6).Value = 1212
$ws.Cells.Item(13, 7).Value = 958
$ws.Cells.Item(13, 8).Value = 32.8
$ws.Cells.Item(13, 9).Value = 25.9

$ws.Cells.Item(14, 5).Value = 39
$ws.Cells.Item(14, 6).Value = 1158
$ws.Cells.Item(14, 7).Value = 1392
$ws.Cells.Item(14, 8).Value = 29.7
$ws.Cells.Item(14, 9).Value = 35.7

# ---------------------------------------------------------------------------
# 4) Meta_ext: bump as_of_utc / build_version.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Meta_ext")
$ws.Cells.Item(2, 2).Value = "2025-10-27T19:00:00Z"
$ws.Cells.Item(2, 4).Value = 17
